# Add a "Distributor_Add" master-data sheet after Sheet1.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Distributor_Add"

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Distributor Code"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Contact Person name"
$ws.Range("D1").Value = "Mobile Number"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Address1"
$ws.Range("G1").Value = "Address2"
$ws.Range("H1").Value = "Address3"
$ws.Range("I1").Value = "Pin Code"
$ws.Range("J1").Value = "Pan Number"
$ws.Range("K1").Value = "FSSAI Licence number"

# --- Sample data row ----------------------------------------------------
$ws.Range("A2").Value = "DB5102"
$ws.Range("B2").Value = "Aniket Enterprices"
$ws.Range("C2").Value = "Aniket sharma"
$ws.Range("D2").Value = 9457863214
$ws.Range("E2").Value = "aniket.jadhav@heerasoftware.com"
$ws.Range("F2").Value = "krudhani sahara appartment"
$ws.Range("G2").Value = "washing center"
$ws.Range("H2").Value = "mumbai"
$ws.Range("I2").Value = 400701
$ws.Range("J2").Value = "KOPGF5479U"
$ws.Range("K2").Value = 10017022012345

# Hyperlink the e-mail address cell (added before the number format below so
# that the Hyperlink cell-style is registered first, matching the authored
# style table order).
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:aniket.jadhav@heerasoftware.com", "", "", "aniket.jadhav@heerasoftware.com") | Out-Null

$ws.Range("K2").NumberFormat = "0.00"

# --- Column widths (best-effort match of the authored layout) ----------
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 18.33
$ws.Columns.Item(4).ColumnWidth = 15.33
$ws.Columns.Item(5).ColumnWidth = 29.33
$ws.Columns.Item(6).ColumnWidth = 23.33
$ws.Columns.Item(7).ColumnWidth = 12.67
$ws.Columns.Item(10).ColumnWidth = 11.67
$ws.Columns.Item(11).ColumnWidth = 18.33

# --- View state: scroll so column C is the first visible one, and
#     leave the active selection on K4, matching the authored sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("K4").Select() | Out-Null
